# Generate Report for Handback
# Update timestamps / status recorded during handback report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 2ad63889... and 5c660467... rows
$wsOverview.Range("G2").Value = "2016-08-27 12:15:20"
$wsOverview.Range("G4").Value = "2016-08-27 12:15:20"

# zh-cn sheet: Priority changes from "ht" to "mt" for 2ad63889... and 5c660467... rows
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# zh-cn sheet: Correspond Handoff/Handback Datetime updates
$wsZhCn.Range("H2").Value = "2016-08-27 12:15:16"
$wsZhCn.Range("H4").Value = "2016-08-27 12:15:16"
$wsZhCn.Range("K2").Value = "2016-08-27 12:15:31"
$wsZhCn.Range("K4").Value = "2016-08-27 12:15:31"

# de-de sheet: Priority changes from "ht" to "mt" for 2ad63889... and 5c660467... rows
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# de-de sheet: Correspond Handback Datetime update
$wsDeDe.Range("K2").Value = "2016-08-27 12:15:37"
$wsDeDe.Range("K4").Value = "2016-08-27 12:15:37"
